# fix: fix bug of main port wire direction
#
# Adds a new port row to the top-level "uart" module port table (sheet1):
#   Port-name = test_temp, InOut = output, Width = 1,
#   Wire-name = tx_busy, Port-comment = test_port

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 16

$ws.Range("A$newRow").Value = "test_temp"
$ws.Range("B$newRow").Value = "output"
$ws.Range("C$newRow").Value = 1
$ws.Range("C$newRow").HorizontalAlignment = -4131
$ws.Range("D$newRow").Value = "tx_busy"
$ws.Range("E$newRow").Value = "test_port"
$ws.Rows.Item($newRow).RowHeight = 16
